$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 9a31e08b-5d28-475a-af66-0dcd1545a378.md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3, Status column B, Latest Handoff Datetime column D
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 06:16:13"

# de-de sheet: row 3, Status column B, Latest Handoff Datetime column D
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 06:16:15"
